$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.958.38"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "3.273.28"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'573.14"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "'177.63"
$ws.Range("E6").Value = "  -4.85%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  +3.69%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -3.25%  "
$ws.Range("D10").Value = "'6.70"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "'0.399"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").Value = "3.845.06"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("E13").Value = "  -3.92%  "
$ws.Range("D14").Value = "65.998.46"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "'26.45"
$ws.Range("E15").Value = "  -3.76%  "
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("D17").Value = "3.271.52"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "'434.66"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("D19").Value = "'5.55"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D20").Value = "'13.14"
$ws.Range("E20").Value = "  -3.37%  "
$ws.Range("D21").Value = "'7.39"
$ws.Range("E21").Value = "  -4.58%  "
$ws.Range("D22").Value = "'72.04"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "3.418.93"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "'0.505"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("D26").Value = "'0.196"
$ws.Range("E26").Value = "  +3.52%  "
$ws.Range("E27").Value = "  -5.64%  "
$ws.Range("D28").Value = "'8.86"
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'1.93"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("D31").Value = "'22.27"
$ws.Range("E31").Value = "  -2.92%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'5.14"
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("D34").Value = "'6.58"
$ws.Range("E34").Value = "  -3.41%  "
$ws.Range("E35").Value = "  -5.48%  "
$ws.Range("D36").Value = "'159.43"
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("E37").Value = "  -5.94%  "
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").Value = "2.755.83"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'0.776"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").Value = "'4.31"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "'6.02"
$ws.Range("E44").Value = "  -3.96%  "
$ws.Range("D45").Value = "'0.0655"
$ws.Range("E45").Value = "  -2.97%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'320.98"
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.28"
$ws.Range("E47").Value = "  -5.75%  "
$ws.Range("D48").Value = "'23.29"
$ws.Range("E48").Value = "  -6.48%  "
$ws.Range("D49").Value = "'0.0266"
$ws.Range("E49").Value = "  -2.90%  "
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.03%  "
